$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 124.44444
$ws.Range("K28").Value = 124.44444
$ws.Range("M28").Value = 360.55556
$ws.Range("H55").Value = 180.57143
$ws.Range("I55").Value = 172.4
$ws.Range("J55").Value = 201
$ws.Range("K55").Value = 172.4
$ws.Range("L55").Value = 201
$ws.Range("M55").Value = 41.59999999999999
$ws.Range("N55").Value = -629
$ws.Range("H98").Value = 4280.077
$ws.Range("I98").Value = 3894.6
$ws.Range("J98").Value = 5565
$ws.Range("K98").Value = 3894.6
$ws.Range("L98").Value = 5565
$ws.Range("M98").Value = -2396.6
$ws.Range("N98").Value = -8561
$ws.Range("H111").Value = 1379.8
$ws.Range("I111").Value = 974.75
$ws.Range("K111").Value = 2924.25
$ws.Range("M111").Value = 142.75
$ws.Range("H122").Value = 4280.077
$ws.Range("I122").Value = 3894.6
$ws.Range("J122").Value = 5565
$ws.Range("K122").Value = 11683.8
$ws.Range("L122").Value = 16695
$ws.Range("M122").Value = -9233.799999999999
$ws.Range("N122").Value = -21595
$ws.Range("H132").Value = 1217.238
$ws.Range("I132").Value = 1119.5
$ws.Range("K132").Value = 3358.5
$ws.Range("M132").Value = -828.5
$ws.Range("H138").Value = 2404.9158
$ws.Range("J138").Value = 2225.5085
$ws.Range("L138").Value = 6676.5255
$ws.Range("N138").Value = -16956.5255
$ws.Range("H139").Value = 72079.664
$ws.Range("J139").Value = 72079.664
$ws.Range("L139").Value = 72079.664
$ws.Range("N139").Value = -82359.664
$ws.Range("H140").Value = 83993.336
$ws.Range("J140").Value = 83993.336
$ws.Range("L140").Value = 83993.336
$ws.Range("N140").Value = -94353.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5425.0625
$ws.Range("I32").Value = 3005.7576
$ws.Range("K32").Value = 3005.7576
$ws.Range("M32").Value = -2718.7576
$ws.Range("H45").Value = 1549.25
$ws.Range("J45").Value = 1699
$ws.Range("L45").Value = 1699
$ws.Range("N45").Value = -2453
$ws.Range("H122").Value = 1323.0834
$ws.Range("I122").Value = 1353.4546
$ws.Range("J122").Value = 989
$ws.Range("K122").Value = 4060.3638
$ws.Range("L122").Value = 2967
$ws.Range("M122").Value = -1610.3638
$ws.Range("N122").Value = -7867
$ws.Range("H132").Value = 2272.3125
$ws.Range("I132").Value = 1971.1333
$ws.Range("K132").Value = 5913.3999
$ws.Range("M132").Value = -3383.3999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1648.3334
$ws.Range("I86").Value = 1630
$ws.Range("J86").Value = 1666.6666
$ws.Range("K86").Value = 1630
$ws.Range("L86").Value = 1666.6666
$ws.Range("M86").Value = -507
$ws.Range("N86").Value = -3912.6666
$ws.Range("H89").Value = 1648.3334
$ws.Range("I89").Value = 1630
$ws.Range("J89").Value = 1666.6666
$ws.Range("K89").Value = 8150
$ws.Range("L89").Value = 8333.333000000001
$ws.Range("M89").Value = -2534
$ws.Range("N89").Value = -19565.333
$ws.Range("H105").Value = 2680.4211
$ws.Range("J105").Value = 4744.5
$ws.Range("L105").Value = 4744.5
$ws.Range("N105").Value = -8238.5
$ws.Range("H107").Value = 722.2857
$ws.Range("I107").Value = 608.6
$ws.Range("K107").Value = 608.6
$ws.Range("M107").Value = 1311.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2419.0908
$ws.Range("I31").Value = 2016.5834
$ws.Range("K31").Value = 2016.5834
$ws.Range("M31").Value = -1721.5834
$ws.Range("H34").Value = 2419.0908
$ws.Range("I34").Value = 2016.5834
$ws.Range("K34").Value = 2016.5834
$ws.Range("M34").Value = -1814.5834
$ws.Range("H94").Value = 1413.72
$ws.Range("I94").Value = 1472.3077
$ws.Range("J94").Value = 1350.25
$ws.Range("K94").Value = 1472.3077
$ws.Range("L94").Value = 1350.25
$ws.Range("M94").Value = -1021.3077
$ws.Range("N94").Value = -2252.25
$ws.Range("H134").Value = 3308.3845
$ws.Range("I134").Value = 3618.3333
$ws.Range("K134").Value = 10854.9999
$ws.Range("M134").Value = -8319.999899999999
$ws.Range("H135").Value = 34693.832
$ws.Range("J135").Value = 34693.832
$ws.Range("L135").Value = 34693.832
$ws.Range("N135").Value = -44833.832

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1999
$ws.Range("J80").Value = 2500
$ws.Range("L80").Value = 7500
$ws.Range("N80").Value = -9372
$ws.Range("H83").Value = 1999
$ws.Range("J83").Value = 2500
$ws.Range("L83").Value = 22500
$ws.Range("N83").Value = -31860
$ws.Range("H112").Value = 51013.75
$ws.Range("I112").Value = 27
$ws.Range("J112").Value = 68009.336
$ws.Range("K112").Value = 81
$ws.Range("L112").Value = 204028.008
$ws.Range("M112").Value = 1027
$ws.Range("N112").Value = -206244.008
$ws.Range("H131").Value = 15801.9
$ws.Range("J131").Value = 16613.053
$ws.Range("L131").Value = 49839.159
$ws.Range("N131").Value = -59919.159
$ws.Range("H140").Value = 2304.0476
$ws.Range("I140").Value = 1773
$ws.Range("K140").Value = 5319
$ws.Range("M140").Value = -139
$ws.Range("H141").Value = 2890.4583
$ws.Range("I141").Value = 2730.3809
$ws.Range("K141").Value = 8191.1427
$ws.Range("M141").Value = -3011.1427

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1875
$ws.Range("I113").Value = 1400
$ws.Range("J113").Value = 1993.75
$ws.Range("K113").Value = 1400
$ws.Range("L113").Value = 1993.75
$ws.Range("M113").Value = 770
$ws.Range("N113").Value = -6333.75
$ws.Range("H122").Value = 1396.3334
$ws.Range("I122").Value = 1316.6666
$ws.Range("J122").Value = 1515.8334
$ws.Range("K122").Value = 3949.9998
$ws.Range("L122").Value = 4547.5002
$ws.Range("M122").Value = -1499.9998
$ws.Range("N122").Value = -9447.5002
$ws.Range("H132").Value = 2266236.2
$ws.Range("I132").Value = 2750430
$ws.Range("K132").Value = 8251290
$ws.Range("M132").Value = -8248760
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2974.3635
$ws.Range("I136").Value = 2968.111
$ws.Range("K136").Value = 8904.332999999999
$ws.Range("M136").Value = -6354.332999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 47000
$ws.Range("J70").Value = 47000
$ws.Range("L70").Value = 47000
$ws.Range("N70").Value = -47630
$ws.Range("H73").Value = 47000
$ws.Range("J73").Value = 47000
$ws.Range("L73").Value = 47000
$ws.Range("N73").Value = -49184
$ws.Range("H126").Value = 1886.5
$ws.Range("I126").Value = 1667.1111
$ws.Range("K126").Value = 5001.3333
$ws.Range("M126").Value = -2531.3333
$ws.Range("H132").Value = 1332.6818
$ws.Range("I132").Value = 1197.8125
$ws.Range("K132").Value = 3593.4375
$ws.Range("M132").Value = -1063.4375

Write-Output "Applied all changes"